# Fruta / hortaliza, semanal
# Update existing rows (2-4), shift old row5 data down into new rows 6-7,
# and set the new values for row5 per the published diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: date + price updates ---
$ws.Range("D2").Value = 44995
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 5500
$ws.Range("P2").Value = 5750
$ws.Range("S2").Value = 2875

# --- Row 3: date, quality, volume/price and origin updates ---
$ws.Range("D3").Value = 45273
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 110
$ws.Range("N3").Value = 9000
$ws.Range("O3").Value = 9500
$ws.Range("P3").Value = 9273
$ws.Range("R3").Value = "Región de Ñuble"
$ws.Range("S3").Value = 4636

# --- Row 4: date, quality, price and origin updates ---
$ws.Range("D4").Value = 45273
$ws.Range("L4").Value = "Segunda"
$ws.Range("N4").Value = 8000
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 8000
$ws.Range("R4").Value = "Región de Ñuble"
$ws.Range("S4").Value = 4000

# --- Row 5: becomes the former row3 "Especial" record (date moved back) ---
$ws.Range("D5").Value = 45008
$ws.Range("L5").Value = "Especial"
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 7000
$ws.Range("O5").Value = 7000
$ws.Range("P5").Value = 7000
$ws.Range("R5").Value = "Provincia de Linares"
$ws.Range("S5").Value = 3500

# --- Row 6 (new): former row4 "Primera" record ---
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 45008
$ws.Range("D6").NumberFormat = $ws.Range("D5").NumberFormat
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100101
$ws.Range("H6").Value = "Berries"
$ws.Range("I6").Value = 100101004
$ws.Range("J6").Value = "Frambuesa"
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 6000
$ws.Range("O6").Value = 6000
$ws.Range("P6").Value = 6000
$ws.Range("Q6").Value = "$/bandeja 2 kilos"
$ws.Range("R6").Value = "Provincia de Linares"
$ws.Range("S6").Value = 3000
$ws.Range("T6").Value = 2

# --- Row 7 (new): former row2 "Primera" record ---
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C7").Value = "Ñuble"
$ws.Range("D7").Value = 44991
$ws.Range("D7").NumberFormat = $ws.Range("D5").NumberFormat
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100101
$ws.Range("H7").Value = "Berries"
$ws.Range("I7").Value = 100101004
$ws.Range("J7").Value = "Frambuesa"
$ws.Range("K7").Value = "Sin especificar"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 6000
$ws.Range("O7").Value = 6000
$ws.Range("P7").Value = 6000
$ws.Range("Q7").Value = "$/bandeja 2 kilos"
$ws.Range("R7").Value = "Provincia de Linares"
$ws.Range("S7").Value = 3000
$ws.Range("T7").Value = 2
